# Renames the embedded-picture labels (wp:docPr/@name and pic:cNvPr/@name)
# on the three Pearson/BTEC logo drawings living in the document's
# headers/footers:
#   - footer (default),  docPr id="1" : image1.png -> image2.png
#   - footer (first page), docPr id="2" : image1.png -> image2.png
#   - header (first page), docPr id="3" : image2.jpg -> image1.jpg
#
# The Word object model has no writable InlineShape.Name, so each drawing
# is replaced in place: capture the exact OOXML of the paragraph that
# hosts it (so no other formatting is disturbed), patch just the `name`
# attribute, delete the shape, and re-insert the patched paragraph at the
# same spot.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoDrawing {
    param(
        [object]$story,     # a Range whose .InlineShapes holds the target drawing (e.g. Footer.Range)
        [int]$shapeIndex,   # 1-based InlineShapes index within that story
        [string]$oldName,   # current wp:docPr / pic:cNvPr name value
        [string]$newName    # replacement name value
    )

    $shp = $story.InlineShapes.Item($shapeIndex)
    $shpRange = $shp.Range

    # Grab the full surrounding paragraph XML (pPr + every run) before
    # touching anything, so formatting/other runs survive untouched.
    $paraXml = $shpRange.Paragraphs.Item(1).Range.WordOpenXML
    if ($paraXml.Length -lt 200) {
        # Some hosts hand back a collapsed range for Paragraphs on a 1-char
        # shape anchor; fall back to expanding the shape's own range.
        $expanded = $shpRange.Duplicate
        $expanded.MoveStart(1, -200) | Out-Null
        $expanded.MoveEnd(1, 200) | Out-Null
        $paraXml = $expanded.Paragraphs.Item(1).Range.WordOpenXML
    }

    $patched = $paraXml.Replace('name="' + $oldName + '"', 'name="' + $newName + '"')

    $shpRange.Delete()
    $shpRange.InsertXML($patched)
}

# --- Footer (default / "odd") page footer: id="1", image1.png -> image2.png
$footerDefault = $sec.Footers.Item(1)
Rename-LogoDrawing -story $footerDefault.Range -shapeIndex 1 -oldName "image1.png" -newName "image2.png"

# --- Footer (first page): id="2", image1.png -> image2.png
$footerFirst = $sec.Footers.Item(2)
Rename-LogoDrawing -story $footerFirst.Range -shapeIndex 1 -oldName "image1.png" -newName "image2.png"

# --- Header (first page): id="3", image2.jpg -> image1.jpg
$headerFirst = $sec.Headers.Item(2)
Rename-LogoDrawing -story $headerFirst.Range -shapeIndex 1 -oldName "image2.jpg" -newName "image1.jpg"

Write-Host "Done renaming logo drawings."
